$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data (D:K) to (E:L)
$ws.Columns("D").Insert()

# Copy formatting (style) from the now-shifted column E into new column D for each data block
$ws.Range("E7:E35").Copy($ws.Range("D7:D35"))
$ws.Range("E38:E77").Copy($ws.Range("D38:D77"))
$ws.Range("E80:E102").Copy($ws.Range("D80:D102"))

# Populate new column D with the latest period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 170756000
$ws.Range("D9").Value = 79419000
$ws.Range("D10").Value = 91337000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -693000
$ws.Range("D15").Value = 28430000
$ws.Range("D17").Value = 143834000
$ws.Range("D18").Value = 26922000
$ws.Range("D20").Value = 5908000
$ws.Range("D21").Value = 65032000
$ws.Range("D22").Value = 7957000
$ws.Range("D23").Value = 24873000
$ws.Range("D24").Value = 5638000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 19235000
$ws.Range("D27").Value = 18652000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 718000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -5908000
$ws.Range("D33").Value = 19370000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 19370000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 5204000
$ws.Range("D43").Value = 29164000
$ws.Range("D44").Value = 2771000
$ws.Range("D45").Value = 14288000
$ws.Range("D46").Value = 51427000
$ws.Range("D47").Value = 6245000
$ws.Range("D48").Value = 131473000
$ws.Range("D49").Value = 310197000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 32522000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 531864000
$ws.Range("D57").Value = 27018000
$ws.Range("D58").Value = 10255000
$ws.Range("D59").Value = 27147000
$ws.Range("D60").Value = 64420000
$ws.Range("D61").Value = 166250000
$ws.Range("D62").Value = 107310000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 347775000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 58753000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 184089000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 19370000
$ws.Range("D83").Value = 32202000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 43602000
$ws.Range("D91").Value = -20758000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -63145000
$ws.Range("D96").Value = -13410000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -25989000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -45532000
